# Frankrig og Tyskland skulle lige bytte index i frac excel
# (France and Germany needed to swap their index columns: C <-> D)
#
# Columns: B=Denmark(?), C=France, D=Germany, E/F=others.
# For rows 2-31 the contents of column C and column D are swapped.
# Rows 2-22 carried a "shared" formula in D (D{r} = D{r+1}/AVERAGE(B{r+1}/B{r},C{r+1}/C{r},E{r+1}/E{r},F{r+1}/F{r}))
# while C held a plain value; after the swap C holds the (column-swapped) formula
# and D holds the plain value that used to live in C. Rows 23-31 only ever held
# plain values in both columns, so those are a straight swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original column C values (rows 2-31) before any writes happen,
# since row r's new D value needs the original C value of that same row.
$origC = @{}
for ($r = 2; $r -le 31; $r++) {
    $origC[$r] = $ws.Cells.Item($r, 3).Value2
}

# Rows 2-22: rewrite the formula (that used to live in D, referencing column C)
# into C, referencing column D instead; put the old plain C value into D.
for ($r = 2; $r -le 22; $r++) {
    $r1 = $r + 1
    $ws.Cells.Item($r, 3).Formula = "=C$r1/AVERAGE(B$r1/B$r,D$r1/D$r,E$r1/E$r,F$r1/F$r)"
    $ws.Cells.Item($r, 4).Value = $origC[$r]
}

# Rows 23-31: plain values in both columns - swap them directly.
for ($r = 23; $r -le 31; $r++) {
    $origD = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $origD
    $ws.Cells.Item($r, 4).Value = $origC[$r]
}

# The author's selection ended up on E6 before saving.
$ws.Range("E6").Select() | Out-Null
